# Update performance dashboard 2025-12-19 00:32
#
# Refreshes the "gemini-3-pro" / Pattern1-Pure Data result row (row 3) on
# both the "Summary" sheet and the "Pattern1-Pure Data" sheet with the
# latest equity / return figures.
#
# The currency-style cells (C/D/E, prefixed with the Yen sign) are written
# as plain text already. The percentage-style cells (F/G/I/J/K) look like
# numbers to Excel's input parser, so they are temporarily forced to the
# Text number format while the literal string is written, then the
# formatting is cleared again so the cells end up as plain, unformatted
# text -- matching how the sheet originally stored these values.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Summary", "Pattern1-Pure Data")
$percentCells = @("F3", "G3", "I3", "J3", "K3")

$newValues = @{
    "C3" = "¥1,000,000.00"
    "D3" = "¥1,001,002.00"
    "E3" = "¥+1,002.00"
    "F3" = "+0.10%"
    "G3" = "+28.71%"
    "I3" = "0.00%"
    "J3" = "100.0%"
    "K3" = "0.1002%"
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Plain-text currency cells -- no special handling needed.
    $ws.Range("C3").Value = $newValues["C3"]
    $ws.Range("D3").Value = $newValues["D3"]
    $ws.Range("E3").Value = $newValues["E3"]

    # Percentage-looking cells: force Text format so the literal string is
    # kept (instead of being parsed into a numeric percentage), then strip
    # the formatting back off so no style is left behind on the cell.
    foreach ($addr in $percentCells) {
        $ws.Range($addr).NumberFormat = "@"
    }
    foreach ($addr in $percentCells) {
        $ws.Range($addr).Value = $newValues[$addr]
    }
    foreach ($addr in $percentCells) {
        $ws.Range($addr).ClearFormats()
    }
}
